$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order timestamps updated) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

$ws1.Name = "GNG_TO-16504777830103743"
$ws2.Name = "NB_TO-16504777856354136"
$ws3.Name = "RS_TO-16504777856363764"
$ws4.Name = "TOL_TO-16504777856833777"
$ws5.Name = "vSAT_TO-16504777857473779"

# --- Sheet 1 (GNG) ---
$ws1.Range("B2").Value = "go_stims-165047778296838.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777829934082.csv"
$ws1.Range("B4").Value = "go_stims-16504777829943786.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777830094068.csv"

# --- Sheet 2 (NB) ---
$ws2.Range("B2").Value = "ZB-match_1-16504777830634098.csv"
$ws2.Range("B3").Value = "TB-16504777843013763.csv"
$ws2.Range("B4").Value = "ZB-match_6-1650477783118411.csv"
$ws2.Range("B5").Value = "OB-16504777834143763.csv"
$ws2.Range("B6").Value = "ZB-match_1-1650477783253409.csv"
$ws2.Range("B7").Value = "OB-16504777836364095.csv"
$ws2.Range("B8").Value = "OB-1650477783847409.csv"
$ws2.Range("B9").Value = "TB-16504777856143782.csv"
$ws2.Range("B10").Value = "TB-1650477784574378.csv"

# --- Sheet 3 (RS) ---
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4.Range("B2").Value = "MM_stims-16504777856513789.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777856383796.csv"
$ws4.Range("B4").Value = "MM_stims-16504777856674047.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777856513789.csv"
$ws4.Range("B6").Value = "MM_stims-16504777856833777.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777856674047.csv"

# --- Sheet 5 (vSAT) ---
$ws5.Range("B2").Value = "SAT_stims-16504777856863787.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504777857153788.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504777857313766.csv"
$ws5.Range("B5").Value = "SAT_stims-16504777856993797.csv"
